$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-24 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-25 Thursday", 2) | Out-Null
$d.Content.Find.Execute("49+42=91", $true, $false, $false, $false, $false, $true, 1, $false, "95-68=27", 2) | Out-Null
$d.Content.Find.Execute("4+49=53", $true, $false, $false, $false, $false, $true, 1, $false, "46-26=20", 2) | Out-Null
$d.Content.Find.Execute("74-44=30", $true, $false, $false, $false, $false, $true, 1, $false, "90-22=68", 2) | Out-Null
$d.Content.Find.Execute("19+79=98", $true, $false, $false, $false, $false, $true, 1, $false, "18+39=57", 2) | Out-Null
$d.Content.Find.Execute("32+9=41", $true, $false, $false, $false, $false, $true, 1, $false, "39+2=41", 2) | Out-Null
$d.Content.Find.Execute("69-42=27", $true, $false, $false, $false, $false, $true, 1, $false, "1+5=6", 2) | Out-Null
$d.Content.Find.Execute("6+88=94", $true, $false, $false, $false, $false, $true, 1, $false, "88-63=25", 2) | Out-Null
$d.Content.Find.Execute("14+55=69", $true, $false, $false, $false, $false, $true, 1, $false, "8+32=40", 2) | Out-Null
$d.Content.Find.Execute("71-70=1", $true, $false, $false, $false, $false, $true, 1, $false, "57+39=96", 2) | Out-Null
$d.Content.Find.Execute("7+81=88", $true, $false, $false, $false, $false, $true, 1, $false, "52-21=31", 2) | Out-Null
$d.Content.Find.Execute("1+67=68", $true, $false, $false, $false, $false, $true, 1, $false, "73+18=91", 2) | Out-Null
$d.Content.Find.Execute("25-24=1", $true, $false, $false, $false, $false, $true, 1, $false, "50+48=98", 2) | Out-Null
$d.Content.Find.Execute("49+3=52", $true, $false, $false, $false, $false, $true, 1, $false, "34+25=59", 2) | Out-Null
$d.Content.Find.Execute("70+25=95", $true, $false, $false, $false, $false, $true, 1, $false, "16+15=31", 2) | Out-Null
$d.Content.Find.Execute("80+9=89", $true, $false, $false, $false, $false, $true, 1, $false, "42+18=60", 2) | Out-Null
$d.Content.Find.Execute("58+4=62", $true, $false, $false, $false, $false, $true, 1, $false, "67+32=99", 2) | Out-Null
$d.Content.Find.Execute("54+7=61", $true, $false, $false, $false, $false, $true, 1, $false, "49-23=26", 2) | Out-Null
$d.Content.Find.Execute("76+23=99", $true, $false, $false, $false, $false, $true, 1, $false, "10-0=10", 2) | Out-Null
$d.Content.Find.Execute("40-17=23", $true, $false, $false, $false, $false, $true, 1, $false, "16+33=49", 2) | Out-Null
$d.Content.Find.Execute("66-36=30", $true, $false, $false, $false, $false, $true, 1, $false, "75-4=71", 2) | Out-Null
$d.Content.Find.Execute("70+11=81", $true, $false, $false, $false, $false, $true, 1, $false, "86-2=84", 2) | Out-Null
$d.Content.Find.Execute("41+30=71", $true, $false, $false, $false, $false, $true, 1, $false, "48+42=90", 2) | Out-Null
$d.Content.Find.Execute("88-22=66", $true, $false, $false, $false, $false, $true, 1, $false, "63+18=81", 2) | Out-Null
$d.Content.Find.Execute("51-6=45", $true, $false, $false, $false, $false, $true, 1, $false, "74+25=99", 2) | Out-Null
$d.Content.Find.Execute("26+59=85", $true, $false, $false, $false, $false, $true, 1, $false, "39-5=34", 2) | Out-Null
$d.Content.Find.Execute("51+8=59", $true, $false, $false, $false, $false, $true, 1, $false, "14+75=89", 2) | Out-Null
$d.Content.Find.Execute("17+75=92", $true, $false, $false, $false, $false, $true, 1, $false, "95-17=78", 2) | Out-Null
$d.Content.Find.Execute("50+24=74", $true, $false, $false, $false, $false, $true, 1, $false, "49-15=34", 2) | Out-Null
$d.Content.Find.Execute("16+50=66", $true, $false, $false, $false, $false, $true, 1, $false, "21+42=63", 2) | Out-Null
$d.Content.Find.Execute("23+22=45", $true, $false, $false, $false, $false, $true, 1, $false, "12+0=12", 2) | Out-Null
$d.Content.Find.Execute("12+57=69", $true, $false, $false, $false, $false, $true, 1, $false, "21+56=77", 2) | Out-Null
$d.Content.Find.Execute("32+11=43", $true, $false, $false, $false, $false, $true, 1, $false, "65+28=93", 2) | Out-Null
$d.Content.Find.Execute("40+40=80", $true, $false, $false, $false, $false, $true, 1, $false, "70-69=1", 2) | Out-Null
$d.Content.Find.Execute("61-21=40", $true, $false, $false, $false, $false, $true, 1, $false, "72+25=97", 2) | Out-Null
$d.Content.Find.Execute("24+27=51", $true, $false, $false, $false, $false, $true, 1, $false, "61+9=70", 2) | Out-Null
$d.Content.Find.Execute("2+46=48", $true, $false, $false, $false, $false, $true, 1, $false, "33+0=33", 2) | Out-Null
$d.Content.Find.Execute("82-9=73", $true, $false, $false, $false, $false, $true, 1, $false, "12+35=47", 2) | Out-Null
$d.Content.Find.Execute("0+53=53", $true, $false, $false, $false, $false, $true, 1, $false, "98-51=47", 2) | Out-Null
$d.Content.Find.Execute("47+3=50", $true, $false, $false, $false, $false, $true, 1, $false, "61+33=94", 2) | Out-Null
$d.Content.Find.Execute("93-29=64", $true, $false, $false, $false, $false, $true, 1, $false, "26-19=7", 2) | Out-Null
$d.Content.Find.Execute("30+21=51", $true, $false, $false, $false, $false, $true, 1, $false, "0+38=38", 2) | Out-Null
$d.Content.Find.Execute("15+22=37", $true, $false, $false, $false, $false, $true, 1, $false, "54-26=28", 2) | Out-Null
$d.Content.Find.Execute("2+13=15", $true, $false, $false, $false, $false, $true, 1, $false, "4+60=64", 2) | Out-Null
$d.Content.Find.Execute("8+36=44", $true, $false, $false, $false, $false, $true, 1, $false, "72-42=30", 2) | Out-Null
$d.Content.Find.Execute("16+81=97", $true, $false, $false, $false, $false, $true, 1, $false, "11+38=49", 2) | Out-Null
$d.Content.Find.Execute("39-33=6", $true, $false, $false, $false, $false, $true, 1, $false, "27-24=3", 2) | Out-Null
$d.Content.Find.Execute("67+8=75", $true, $false, $false, $false, $false, $true, 1, $false, "49-20=29", 2) | Out-Null
$d.Content.Find.Execute("57+32=89", $true, $false, $false, $false, $false, $true, 1, $false, "39+42=81", 2) | Out-Null
$d.Content.Find.Execute("70+20=90", $true, $false, $false, $false, $false, $true, 1, $false, "94-59=35", 2) | Out-Null
$d.Content.Find.Execute("22+34=56", $true, $false, $false, $false, $false, $true, 1, $false, "10+83=93", 2) | Out-Null
$d.Content.Find.Execute("47-4=43", $true, $false, $false, $false, $false, $true, 1, $false, "1+44=45", 2) | Out-Null
$d.Content.Find.Execute("26+50=76", $true, $false, $false, $false, $false, $true, 1, $false, "64-3=61", 2) | Out-Null
$d.Content.Find.Execute("93-48=45", $true, $false, $false, $false, $false, $true, 1, $false, "42+14=56", 2) | Out-Null
$d.Content.Find.Execute("8+45=53", $true, $false, $false, $false, $false, $true, 1, $false, "93-87=6", 2) | Out-Null
$d.Content.Find.Execute("15+27=42", $true, $false, $false, $false, $false, $true, 1, $false, "96-28=68", 2) | Out-Null
$d.Content.Find.Execute("44+39=83", $true, $false, $false, $false, $false, $true, 1, $false, "37+40=77", 2) | Out-Null
$d.Content.Find.Execute("89-42=47", $true, $false, $false, $false, $false, $true, 1, $false, "49-43=6", 2) | Out-Null
$d.Content.Find.Execute("69-58=11", $true, $false, $false, $false, $false, $true, 1, $false, "1+8=9", 2) | Out-Null
$d.Content.Find.Execute("97-12=85", $true, $false, $false, $false, $false, $true, 1, $false, "10+69=79", 2) | Out-Null
$d.Content.Find.Execute("11+45=56", $true, $false, $false, $false, $false, $true, 1, $false, "35+63=98", 2) | Out-Null
$d.Content.Find.Execute("23+43=66", $true, $false, $false, $false, $false, $true, 1, $false, "23+38=61", 2) | Out-Null
$d.Content.Find.Execute("82+8=90", $true, $false, $false, $false, $false, $true, 1, $false, "90-6=84", 2) | Out-Null
$d.Content.Find.Execute("39-9=30", $true, $false, $false, $false, $false, $true, 1, $false, "52-50=2", 2) | Out-Null
$d.Content.Find.Execute("32-25=7", $true, $false, $false, $false, $false, $true, 1, $false, "31+13=44", 2) | Out-Null
$d.Content.Find.Execute("72-53=19", $true, $false, $false, $false, $false, $true, 1, $false, "61-3=58", 2) | Out-Null
$d.Content.Find.Execute("69-32=37", $true, $false, $false, $false, $false, $true, 1, $false, "88-76=12", 2) | Out-Null
$d.Content.Find.Execute("29-5=24", $true, $false, $false, $false, $false, $true, 1, $false, "54-2=52", 2) | Out-Null
$d.Content.Find.Execute("95+1=96", $true, $false, $false, $false, $false, $true, 1, $false, "80-60=20", 2) | Out-Null
$d.Content.Find.Execute("37-12=25", $true, $false, $false, $false, $false, $true, 1, $false, "95-94=1", 2) | Out-Null
$d.Content.Find.Execute("72-19=53", $true, $false, $false, $false, $false, $true, 1, $false, "80-5=75", 2) | Out-Null
$d.Content.Find.Execute("58-7=51", $true, $false, $false, $false, $false, $true, 1, $false, "55-41=14", 2) | Out-Null
$d.Content.Find.Execute("66-46=20", $true, $false, $false, $false, $false, $true, 1, $false, "68+2=70", 2) | Out-Null
$d.Content.Find.Execute("27+11=38", $true, $false, $false, $false, $false, $true, 1, $false, "60-50=10", 2) | Out-Null
$d.Content.Find.Execute("39+60=99", $true, $false, $false, $false, $false, $true, 1, $false, "7+54=61", 2) | Out-Null
$d.Content.Find.Execute("16+42=58", $true, $false, $false, $false, $false, $true, 1, $false, "44-21=23", 2) | Out-Null
$d.Content.Find.Execute("82-12=70", $true, $false, $false, $false, $false, $true, 1, $false, "25-5=20", 2) | Out-Null
$d.Content.Find.Execute("63+17=80", $true, $false, $false, $false, $false, $true, 1, $false, "53+35=88", 2) | Out-Null
$d.Content.Find.Execute("63+11=74", $true, $false, $false, $false, $false, $true, 1, $false, "55+26=81", 2) | Out-Null
$d.Content.Find.Execute("38-31=7", $true, $false, $false, $false, $false, $true, 1, $false, "46-11=35", 2) | Out-Null
$d.Content.Find.Execute("27+68=95", $true, $false, $false, $false, $false, $true, 1, $false, "67+14=81", 2) | Out-Null
$d.Content.Find.Execute("43+6=49", $true, $false, $false, $false, $false, $true, 1, $false, "41+8=49", 2) | Out-Null
$d.Content.Find.Execute("13+8=21", $true, $false, $false, $false, $false, $true, 1, $false, "32-15=17", 2) | Out-Null
$d.Content.Find.Execute("46+47=93", $true, $false, $false, $false, $false, $true, 1, $false, "19+77=96", 2) | Out-Null
$d.Content.Find.Execute("77+9=86", $true, $false, $false, $false, $false, $true, 1, $false, "25-3=22", 2) | Out-Null
$d.Content.Find.Execute("24+53=77", $true, $false, $false, $false, $false, $true, 1, $false, "34-29=5", 2) | Out-Null
$d.Content.Find.Execute("64+31=95", $true, $false, $false, $false, $false, $true, 1, $false, "79+4=83", 2) | Out-Null
$d.Content.Find.Execute("41-21=20", $true, $false, $false, $false, $false, $true, 1, $false, "8+58=66", 2) | Out-Null
$d.Content.Find.Execute("90-66=24", $true, $false, $false, $false, $false, $true, 1, $false, "29+47=76", 2) | Out-Null
$d.Content.Find.Execute("28-1=27", $true, $false, $false, $false, $false, $true, 1, $false, "92-54=38", 2) | Out-Null
$d.Content.Find.Execute("16-12=4", $true, $false, $false, $false, $false, $true, 1, $false, "20-6=14", 2) | Out-Null
$d.Content.Find.Execute("51+37=88", $true, $false, $false, $false, $false, $true, 1, $false, "42-32=10", 2) | Out-Null
$d.Content.Find.Execute("38+47=85", $true, $false, $false, $false, $false, $true, 1, $false, "2+53=55", 2) | Out-Null
$d.Content.Find.Execute("79+15=94", $true, $false, $false, $false, $false, $true, 1, $false, "75-41=34", 2) | Out-Null
$d.Content.Find.Execute("49+29=78", $true, $false, $false, $false, $false, $true, 1, $false, "8+66=74", 2) | Out-Null
$d.Content.Find.Execute("23+71=94", $true, $false, $false, $false, $false, $true, 1, $false, "93-2=91", 2) | Out-Null
$d.Content.Find.Execute("38-10=28", $true, $false, $false, $false, $false, $true, 1, $false, "1+56=57", 2) | Out-Null
$d.Content.Find.Execute("51-43=8", $true, $false, $false, $false, $false, $true, 1, $false, "14+32=46", 2) | Out-Null
$d.Content.Find.Execute("5+58=63", $true, $false, $false, $false, $false, $true, 1, $false, "15+57=72", 2) | Out-Null
$d.Content.Find.Execute("65+7=72", $true, $false, $false, $false, $false, $true, 1, $false, "23+67=90", 2) | Out-Null
$d.Content.Find.Execute("93-49=44", $true, $false, $false, $false, $false, $true, 1, $false, "52-38=14", 2) | Out-Null
